# Auto-generated edit script: updates cryptos Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.664.27"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.460.06"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'559.19"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "'161.68"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "68.570.07"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").Value = "'23.49"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "2.458.01"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D19").Value = "'334.01"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").Value = "'66.46"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").Value = "'8.14"
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("D27").Value = "0.0₃0813"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("D28").Value = "'7.17"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D30").Value = "'428.79"
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("E31").Value = "  -4.28%  "
$ws.Range("E32").Value = "  -4.38%  "
$ws.Range("D33").Value = "'158.86"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("D34").Value = "'19.02"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("E40").Value = "  -4.89%  "
$ws.Range("E41").Value = "  -5.93%  "
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").Value = "'129.73"
$ws.Range("E44").Value = "  -4.13%  "
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "'0.559"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").Value = "'0.0908"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("E51").Value = "  -8.58%  "
